# Hangman_v4: add Sheet2 (flow-chart notes for the while-loop / menu logic)
# plus its drawing shapes, and move Sheet1's selection.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new worksheet right after Sheet1 --------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Activate()

# --- Populate cells in the same order the original author typed them ---
# (this keeps the shared-string table indices lined up with the diff)
$ws2.Range("D8").Value = "int userChoice = "
$ws2.Range("D13").Value = "Case 1"
$ws2.Range("H14").Value = "eletek"
$ws2.Range("I14").Value = 2
$ws2.Range("D16").Value = "Case 2"
$ws2.Range("H17").Value = "hos"
$ws2.Range("I18").Value = 4
$ws2.Range("D20").Value = "Cas 3"
$ws2.Range("D23").Value = "userArray"
$ws2.Range("D24").Value = "initalArray"
$ws2.Range("K4").Value = "Array Nehez"
$ws2.Range("M4").Value = "Array kozepes"
$ws2.Range("O4").Value = "Array konntu"
$ws2.Range("H13").Value = "u ArrayList<String> initialWordArray = createInitialWordArray(solutionArrayHard);"
$ws2.Range("H16").Value = "u ArrayList<String> initialWordArray = createInitialWordArray(solutionArrayMEdium)"

# --- Flow-chart shapes for the while-loop sketch ------------------------
$shp1 = $ws2.Shapes.AddShape(1, 239.25, 47, 84, 38.75)
$shp1.Name = "Rectangle 1"
$shp1.TextFrame.Characters().Text = "Show Menu"

$shp2 = $ws2.Shapes.AddShape(1, 233.75, 99.25, 84, 53.5)
$shp2.Name = "Rectangle 3"
$shp2.TextFrame.Characters().Text = "getUserInput()" + [char]10 + "Scanner" + [char]10 + "return "

$shp3 = $ws2.Shapes.AddShape(4, 208, 176.75, 75, 112.75)
$shp3.Name = "Diamond 5"

[void]$ws2.Range("L16").Select()

# --- Restore Sheet1 as the active tab/selection -------------------------
[void]$ws1.Activate()
[void]$ws1.Range("T5").Select()
